$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 453
$ws.Range("D3").Value = 82.5

$ws.Range("C4").Value = 1110
$ws.Range("D4").Value = 92

$ws.Range("C5").Value = 7240

$ws.Range("C6").Value = 70
$ws.Range("D6").Value = 97.40000000000001

$ws.Range("D7").Value = 97.90000000000001

$ws.Range("C8").Value = 8950
